$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 69, shifting existing rows 69-140 down to 70-141
$ws.Rows.Item(69).Insert()

# Populate the newly inserted row 69 with the new record
$ws.Range("A69").Value = 10
$ws.Range("B69").Value = 'Vega Modelo de Temuco'
$ws.Range("C69").Value = 'La Araucanía'
$ws.Range("D69").Value = 45159
$ws.Range("D69").NumberFormat = $ws.Range("D70").NumberFormat
$ws.Range("E69").Value = 9
$ws.Range("F69").Value = 100112010
$ws.Range("G69").Value = 'Achicoria'
$ws.Range("H69").Value = 'Sin especificar'
$ws.Range("I69").Value = 'Primera'
$ws.Range("J69").Value = 65
$ws.Range("K69").Value = 10000
$ws.Range("L69").Value = 10000
$ws.Range("M69").Value = 10000
$ws.Range("N69").Value = '$/caja 18 unidades'
$ws.Range("O69").Value = 'Región Metropolitana'
$ws.Range("P69").Value = 556
$ws.Range("Q69").Value = 18
$ws.Range("R69").Value = 'Hortaliza'
